{"js": "// Update the worksheet date and the 25 division problems.\n// The date lives in the first paragraph of the document body.\nconst body = context.document.body;\nconst dateHits = body.search(\"2024-03-16 Saturday\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\ndateHits.items[0].insertText(\"2024-03-17 Sunday\", \"Replace\");\n\n// The division problems live in the single table, 5 \"content\" rows\n// (table row indices 0, 4, 8, 12, 16 \u2014 the intervening rows are blank\n// spacer rows) x 5 columns each.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst newValues = {\n  0: [\"61\u00f72=\", \"94\u00f76=\", \"57\u00f74=\", \"98\u00f74=\", \"85\u00f72=\"],\n  4: [\"96\u00f78=\", \"19\u00f76=\", \"72\u00f78=\", \"72\u00f73=\", \"18\u00f74=\"],\n  8: [\"22\u00f77=\", \"53\u00f73=\", \"37\u00f75=\", \"94\u00f78=\", \"54\u00f72=\"],\n  12: [\"54\u00f78=\", \"20\u00f76=\", \"22\u00f78=\", \"35\u00f77=\", \"57\u00f74=\"],\n  16: [\"28\u00f79=\", \"98\u00f79=\", \"99\u00f75=\", \"23\u00f76=\", \"53\u00f76=\"],\n};\n\nfor (const rowIndex of Object.keys(newValues)) {\n  const r = Number(rowIndex);\n  const rowVals = newValues[rowIndex];\n  for (let c = 0; c < rowVals.length; c++) {\n    table.getCell(r, c).value = rowVals[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and the 25 division problems.\n$d = $word.ActiveDocument\n\n# The date paragraph is the very first paragraph in the document.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Execute(\"2024-03-16 Saturday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2024-03-17 Sunday\", 2) | Out-Null\n\n# The division problems live in the single table, 5 \"content\" rows\n# (1-based table row numbers 1, 5, 9, 13, 17 -- the intervening rows are\n# blank spacer rows) x 5 columns each.\n$t = $d.Tables(1)\n\n$newValues = @{\n    1  = @(\"61\u00f72=\", \"94\u00f76=\", \"57\u00f74=\", \"98\u00f74=\", \"85\u00f72=\")\n    5  = @(\"96\u00f78=\", \"19\u00f76=\", \"72\u00f78=\", \"72\u00f73=\", \"18\u00f74=\")\n    9  = @(\"22\u00f77=\", \"53\u00f73=\", \"37\u00f75=\", \"94\u00f78=\", \"54\u00f72=\")\n    13 = @(\"54\u00f78=\", \"20\u00f76=\", \"22\u00f78=\", \"35\u00f77=\", \"57\u00f74=\")\n    17 = @(\"28\u00f79=\", \"98\u00f79=\", \"99\u00f75=\", \"23\u00f76=\", \"53\u00f76=\")\n}\n\nforeach ($row in $newValues.Keys) {\n    $rowVals = $newValues[$row]\n    for ($c = 1; $c -le $rowVals.Length; $c++) {\n        $t.Cell($row, $c).Range.Text = $rowVals[$c - 1]\n    }\n}\n"}
